# Applies the cryptos price/volume refresh described in the commit:
# "Updated cryptos list on Thu Jun  1 10:11:14 UTC 2023 with GitHub Actions"
#
# Columns B (Coin) and C (Link) occasionally swap two adjacent rows'
# contents (ranking re-sort); columns D (Price) and E (Volume(1h)) get
# refreshed text values. Every written value is forced through a
# temporary "@" (Text) number format so Excel's COM layer does not
# silently reinterpret numeric-looking strings (e.g. "1.000", "20.70")
# as actual numbers -- the source data are plain text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value="26.906.68"},
    @{Cell="E2"; Value="  -0.90%  "},
    @{Cell="D3"; Value="1.861.83"},
    @{Cell="E3"; Value="  -0.53%  "},
    @{Cell="D4"; Value="1.000"},
    @{Cell="E4"; Value="  -0.03%  "},
    @{Cell="D5"; Value="304.97"},
    @{Cell="E5"; Value="  -0.79%  "},
    @{Cell="D6"; Value="1.000"},
    @{Cell="D7"; Value="0.5067"},
    @{Cell="E7"; Value="  +0.32%  "},
    @{Cell="D8"; Value="0.3619"},
    @{Cell="E8"; Value="  -3.51%  "},
    @{Cell="D9"; Value="0.07172"},
    @{Cell="E9"; Value="  +0.19%  "},
    @{Cell="D10"; Value="0.8954"},
    @{Cell="E10"; Value="  +0.61%  "},
    @{Cell="D11"; Value="20.70"},
    @{Cell="E11"; Value="  -0.08%  "},
    @{Cell="B12"; Value="TRON"},
    @{Cell="C12"; Value="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"},
    @{Cell="D12"; Value="0.07444"},
    @{Cell="E12"; Value="  -1.66%  "},
    @{Cell="B13"; Value="WrappedEther"},
    @{Cell="C13"; Value="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"},
    @{Cell="D13"; Value="1.834.11"},
    @{Cell="E13"; Value="  -2.05%  "},
    @{Cell="D14"; Value="92.95"},
    @{Cell="E14"; Value="  +3.96%  "},
    @{Cell="D15"; Value="5.235"},
    @{Cell="E15"; Value="  -1.71%  "},
    @{Cell="D16"; Value="1.001"},
    @{Cell="E16"; Value="  -0.05%  "},
    @{Cell="D17"; Value="0.000008477"},
    @{Cell="E17"; Value="  -0.30%  "},
    @{Cell="E18"; Value="  +0.06%  "},
    @{Cell="E19"; Value="  -0.09%  "},
    @{Cell="D20"; Value="26.948.91"},
    @{Cell="E20"; Value="  -0.94%  "},
    @{Cell="D21"; Value="5.032"},
    @{Cell="E21"; Value="  -0.96%  "},
    @{Cell="D22"; Value="2.091.41"},
    @{Cell="D23"; Value="10.32"},
    @{Cell="E23"; Value="  -2.71%  "},
    @{Cell="D24"; Value="6.418"},
    @{Cell="E24"; Value="  -1.17%  "},
    @{Cell="D25"; Value="147.93"},
    @{Cell="E25"; Value="  -1.95%  "},
    @{Cell="D26"; Value="1.795"},
    @{Cell="E26"; Value="  -2.60%  "},
    @{Cell="D27"; Value="17.86"},
    @{Cell="E27"; Value="  -0.90%  "},
    @{Cell="D28"; Value="2.060"},
    @{Cell="E28"; Value="  -1.57%  "},
    @{Cell="D29"; Value="113.15"},
    @{Cell="E29"; Value="  +0.19%  "},
    @{Cell="D30"; Value="4.670"},
    @{Cell="E30"; Value="  -2.04%  "},
    @{Cell="D31"; Value="4.675"},
    @{Cell="E31"; Value="  -0.29%  "},
    @{Cell="D32"; Value="0.09254"},
    @{Cell="E32"; Value="  +2.77%  "},
    @{Cell="D33"; Value="0.05077"},
    @{Cell="E33"; Value="  -1.03%  "},
    @{Cell="B34"; Value="HuobiToken"},
    @{Cell="C34"; Value="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"},
    @{Cell="D34"; Value="2.992"},
    @{Cell="E34"; Value="  -3.14%  "},
    @{Cell="B35"; Value="ImmutableX"},
    @{Cell="C35"; Value="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"},
    @{Cell="D35"; Value="0.7457"},
    @{Cell="E35"; Value="  +0.29%  "},
    @{Cell="D36"; Value="1.150"},
    @{Cell="E36"; Value="  -0.95%  "},
    @{Cell="D37"; Value="3.288"},
    @{Cell="E37"; Value="  +7.91%  "},
    @{Cell="D38"; Value="2.516"},
    @{Cell="E38"; Value="  -1.03%  "},
    @{Cell="D39"; Value="0.02000"},
    @{Cell="E39"; Value="  -1.77%  "},
    @{Cell="B40"; Value="TheSandbox"},
    @{Cell="C40"; Value="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"},
    @{Cell="D40"; Value="0.5515"},
    @{Cell="E40"; Value="  +2.60%  "},
    @{Cell="B41"; Value="TrustWalletToken"},
    @{Cell="C41"; Value="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"},
    @{Cell="D41"; Value="1.075"},
    @{Cell="E41"; Value="  -0.11%  "},
    @{Cell="D42"; Value="118.10"},
    @{Cell="E42"; Value="  +2.64%  "},
    @{Cell="D43"; Value="6.483"},
    @{Cell="E43"; Value="  -1.98%  "},
    @{Cell="D44"; Value="8.539"},
    @{Cell="E44"; Value="  +1.22%  "},
    @{Cell="E45"; Value="  -0.71%  "},
    @{Cell="D46"; Value="0.4699"},
    @{Cell="E46"; Value="  +1.03%  "},
    @{Cell="D47"; Value="0.9997"},
    @{Cell="E47"; Value="  -0.03%  "},
    @{Cell="D48"; Value="10.06"},
    @{Cell="E48"; Value="  +0.02%  "},
    @{Cell="D49"; Value="1.563"},
    @{Cell="E49"; Value="  -0.70%  "},
    @{Cell="D50"; Value="37.03"},
    @{Cell="E50"; Value="  +1.30%  "},
    @{Cell="E51"; Value="  -2.61%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $savedFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.NumberFormat = $savedFormat
}
